# Regeneration of en/fr/es country data books after child program paras updated
$wb = $excel.ActiveWorkbook

# --- Update "Programmes pour les enfants" sheet (child program paras) data ---
$ws = $wb.Worksheets.Item("Programmes pour les enfants")
$ws.Activate()
    $ws.Range("F2").Value = 0.39473684210526322
    $ws.Range("G2").Value = 0.39473684210526322
    $ws.Range("H2").Value = 0.39473684210526322
    $ws.Range("F3").Value = 0.30769230769230765
    $ws.Range("G3").Value = 0.30769230769230765
    $ws.Range("H3").Value = 0.30769230769230765
    $ws.Range("F18").Value = 0.7
    $ws.Range("F20").Value = 0.84
    $ws.Range("D21").Value = 0.28260869565217389
    $ws.Range("F21").Value = 0
    $ws.Range("F22").Value = 0
    $ws.Range("D23").Value = 0.28260869565217389
    $ws.Range("F23").Value = 0
    $ws.Range("F24").Value = 0
    $ws.Range("D25").Value = 0.28260869565217389
    $ws.Range("F25").Value = 0
    $ws.Range("F26").Value = 0
    $ws.Range("F27").Value = 1
    $ws.Range("F28").Value = 0
    $ws.Range("F29").Value = 0
    $ws.Range("F30").Value = 1
    $ws.Range("F31").Value = 0
    $ws.Range("F32").Value = 0
    $ws.Range("F33").Value = 1
    $ws.Range("F34").Value = 0
    $ws.Range("F35").Value = 0
    $ws.Range("F36").Value = 1
    $ws.Range("F37").Value = 0
    $ws.Range("F38").Value = 0
    $ws.Range("F39").Value = 1
    $ws.Range("F40").Value = 0
    $ws.Range("F41").Value = 0
    $ws.Range("F42").Value = 0.3
    $ws.Range("F43").Value = 0.5
    $ws.Range("F44").Value = 0.65
    $ws.Range("F45").Value = 0.3
    $ws.Range("F46").Value = 0.49
    $ws.Range("F47").Value = 0.52
    $ws.Range("F48").Value = 0.88
    $ws.Range("D49").Value = 0.78409090909090906
    $ws.Range("E49").Value = 0.78409090909090906
    $ws.Range("F49").Value = 0.78409090909090906
    $ws.Range("G49").Value = 0.78409090909090906
    $ws.Range("H49").Value = 0.78409090909090906
    $ws.Range("D50").Value = 0.88372093023255816
    $ws.Range("E50").Value = 0.88372093023255816
    $ws.Range("F50").Value = 0.88372093023255816
    $ws.Range("G50").Value = 0.88372093023255816
    $ws.Range("H50").Value = 0.88372093023255816
    $ws.Range("F51").Value = 0.86
    $ws.Range("F52").Value = 0
    $ws.Range("F53").Value = 0

# Update the sheet's selection/view to match the latest editing focus
$ws.Range("D2:H53").Select()

# --- Switch the active/selected tab to "Programmes-anémie" ---
$wsAnemie = $wb.Worksheets.Item("Programmes-anémie")
$wsAnemie.Activate()
$wsAnemie.Range("C49").Select()
